$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 was stored as text "79174445" - convert it to a real number
$ws.Range("A2").Value = 79174445

# Row 3: new payment record for phone 79174445 (Cash, 2025-08-23T09:26:33)
# Keep the phone number as text (leading apostrophe forces text entry),
# then reset the cell style back to Normal so no extra number-format style
# gets attached to the cell.
$ws.Range("A3").Value = "'79174445"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 25
$ws.Range("G3").Value = "Cash"
$ws.Range("H3").Value = "2025-08-23T09:26:33"
